# 43_scenecat_block_order.xlsx - reorder block columns from
# (kitchens_1, bedrooms_1, kitchens_2, living_rooms_1, living_rooms_2, bedrooms_2)
# to (living_rooms_1, living_rooms_2, bedrooms_1, bedrooms_2, kitchens_1, kitchens_2)
# and move each row's single "1" marker to the column matching its new position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header order (row 1)
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "living_rooms_2"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("D1").Value = "bedrooms_2"
$ws.Range("E1").Value = "kitchens_1"
$ws.Range("F1").Value = "kitchens_2"

# Rebuild the 0/1 indicator matrix (rows 2-7) for the new column order
$data = @(
    @(0, 0, 0, 1, 0, 0),
    @(0, 0, 0, 0, 1, 0),
    @(1, 0, 0, 0, 0, 0),
    @(0, 0, 1, 0, 0, 0),
    @(0, 0, 0, 0, 0, 1),
    @(0, 1, 0, 0, 0, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowVals[$j]
    }
}
